$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update values
$ws.Range("C2").Value = 11
$ws.Range("C3").Value = 10
$ws.Range("B4").Value = 0.7

# Delete row 6 entirely (shifts dimension from A1:C6 to A1:C5)
$ws.Rows.Item(6).Delete()

# Update column widths: A and C get explicit custom widths (B reverts to default)
$ws.Columns.Item(1).ColumnWidth = 26.285714285714285
$ws.Columns.Item(3).ColumnWidth = 26.571428571428573

# Update selection to C3
$ws.Range("C3").Select()

$wb.Save()
